# Scheduled runner update: refresh Leve profit-calc sheets with latest
# market-board price snapshots (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1965.7192
$ws.Range("I15").Value = 1965.7192
$ws.Range("K15").Value = 5897.1576
$ws.Range("M15").Value = -5728.1576
$ws.Range("H33").Value = 281.5435
$ws.Range("I33").Value = 252.92683
$ws.Range("K33").Value = 252.92683
$ws.Range("M33").Value = -23.92683
$ws.Range("H132").Value = 34918.547
$ws.Range("I132").Value = 34918.547
$ws.Range("K132").Value = 104755.641
$ws.Range("M132").Value = -102225.641
$ws.Range("I137").Value = 2135.2334
$ws.Range("J137").Value = 1860.1578
$ws.Range("K137").Value = 6405.7002
$ws.Range("L137").Value = 5580.4734
$ws.Range("M137").Value = -3855.7002
$ws.Range("N137").Value = -10680.4734

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3359.6938
$ws.Range("J32").Value = 9001.75
$ws.Range("L32").Value = 9001.75
$ws.Range("N32").Value = -9575.75
$ws.Range("H45").Value = 2709.2
$ws.Range("I45").Value = 2428
$ws.Range("J45").Value = 3030.5715
$ws.Range("K45").Value = 2428
$ws.Range("L45").Value = 3030.5715
$ws.Range("M45").Value = -2051
$ws.Range("N45").Value = -3784.5715
$ws.Range("H125").Value = 36000
$ws.Range("J125").Value = 36000
$ws.Range("L125").Value = 36000
$ws.Range("N125").Value = -45840
$ws.Range("H132").Value = 66313.75
$ws.Range("I132").Value = 4102.4
$ws.Range("J132").Value = 169999.33
$ws.Range("K132").Value = 12307.2
$ws.Range("L132").Value = 509997.99
$ws.Range("M132").Value = -9777.199999999999
$ws.Range("N132").Value = -515057.99

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 19270.285
$ws.Range("J81").Value = 19270.285
$ws.Range("L81").Value = 19270.285
$ws.Range("N81").Value = -21392.285
$ws.Range("H84").Value = 19270.285
$ws.Range("J84").Value = 19270.285
$ws.Range("L84").Value = 57810.855
$ws.Range("N84").Value = -68418.855
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13105.8
$ws.Range("I31").Value = 22630.785
$ws.Range("J31").Value = 4771.4375
$ws.Range("K31").Value = 22630.785
$ws.Range("L31").Value = 4771.4375
$ws.Range("M31").Value = -22335.785
$ws.Range("N31").Value = -5361.4375
$ws.Range("H34").Value = 13105.8
$ws.Range("I34").Value = 22630.785
$ws.Range("J34").Value = 4771.4375
$ws.Range("K34").Value = 22630.785
$ws.Range("L34").Value = 4771.4375
$ws.Range("M34").Value = -22428.785
$ws.Range("N34").Value = -5175.4375
$ws.Range("H88").Value = 31000
$ws.Range("J88").Value = 31000
$ws.Range("L88").Value = 31000
$ws.Range("N88").Value = -31812
$ws.Range("H91").Value = 31000
$ws.Range("J91").Value = 31000
$ws.Range("L91").Value = 31000
$ws.Range("N91").Value = -33808
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H132").Value = 18167
$ws.Range("I132").Value = 30530
$ws.Range("J132").Value = 4155.6
$ws.Range("K132").Value = 91590
$ws.Range("L132").Value = 12466.8
$ws.Range("M132").Value = -89060
$ws.Range("N132").Value = -17526.8
$ws.Range("H134").Value = 1027
$ws.Range("I134").Value = 1029.4546
$ws.Range("K134").Value = 3088.3638
$ws.Range("M134").Value = -553.3638000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 762.7
$ws.Range("I5").Value = 716.75
$ws.Range("K5").Value = 2150.25
$ws.Range("M5").Value = -2038.25
$ws.Range("H68").Value = 1330.7576
$ws.Range("J68").Value = 1653.3043
$ws.Range("L68").Value = 4959.9129
$ws.Range("N68").Value = -6581.9129
$ws.Range("H71").Value = 1330.7576
$ws.Range("J71").Value = 1653.3043
$ws.Range("L71").Value = 14879.7387
$ws.Range("N71").Value = -22991.7387
$ws.Range("H75").Value = 5500
$ws.Range("J75").Value = 5500
$ws.Range("L75").Value = 16500
$ws.Range("N75").Value = -18496
$ws.Range("H78").Value = 5500
$ws.Range("J78").Value = 5500
$ws.Range("L78").Value = 49500
$ws.Range("N78").Value = -59484
$ws.Range("H96").Value = 3500
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -16118
$ws.Range("H102").Value = 5720
$ws.Range("J102").Value = 5720
$ws.Range("L102").Value = 17160
$ws.Range("N102").Value = -22028
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").ClearContents()
$ws.Range("I132").Value = 897.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8080.2
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 762.7
$ws.Range("I135").Value = 716.75
$ws.Range("K135").Value = 6450.75
$ws.Range("M135").Value = -3915.75

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 57.833332
$ws.Range("I2").Value = 51.9
$ws.Range("J2").Value = 87.5
$ws.Range("K2").Value = 51.9
$ws.Range("L2").Value = 87.5
$ws.Range("M2").Value = 61.1
$ws.Range("N2").Value = -313.5
$ws.Range("H102").Value = 2332.2856
$ws.Range("I102").Value = 2502.4
$ws.Range("J102").Value = 1907
$ws.Range("K102").Value = 2502.4
$ws.Range("L102").Value = 1907
$ws.Range("M102").Value = -880.4000000000001
$ws.Range("N102").Value = -5151
$ws.Range("H126").Value = 6515.125
$ws.Range("I126").Value = 7108.3335
$ws.Range("J126").Value = 5921.9165
$ws.Range("K126").Value = 21325.0005
$ws.Range("L126").Value = 17765.7495
$ws.Range("M126").Value = -18855.0005
$ws.Range("N126").Value = -22705.7495

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5840
$ws.Range("I68").Value = 3066.6667
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 3066.6667
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -2317.6667
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 5840
$ws.Range("I71").Value = 3066.6667
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 15333.3335
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -11589.3335
$ws.Range("N71").Value = -57488
$ws.Range("H132").Value = 2857
$ws.Range("I132").Value = 2082.5
$ws.Range("J132").Value = 3437.875
$ws.Range("K132").Value = 6247.5
$ws.Range("L132").Value = 10313.625
$ws.Range("M132").Value = -3717.5
$ws.Range("N132").Value = -15373.625

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10226
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H70").Value = 21500
$ws.Range("I70").Value = 20666.666
$ws.Range("J70").Value = 24000
$ws.Range("K70").Value = 20666.666
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = -20351.666
$ws.Range("N70").Value = -24630
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H73").Value = 21500
$ws.Range("I73").Value = 20666.666
$ws.Range("J73").Value = 24000
$ws.Range("K73").Value = 20666.666
$ws.Range("L73").Value = 24000
$ws.Range("M73").Value = -19574.666
$ws.Range("N73").Value = -26184
$ws.Range("H122").Value = 1793
$ws.Range("I122").Value = 1625.8889
$ws.Range("K122").Value = 4877.6667
$ws.Range("M122").Value = -2427.6667
$ws.Range("H132").Value = 2106.36
$ws.Range("I132").Value = 1874.2632
$ws.Range("J132").Value = 2841.3333
$ws.Range("K132").Value = 5622.7896
$ws.Range("L132").Value = 8523.999899999999
$ws.Range("M132").Value = -3092.7896
$ws.Range("N132").Value = -13583.9999

